# Agregue una pregunta mas.
# Adds a new "to do" item (question #18) to the "General" sheet's backlog
# table, and updates the saved selections on both sheets to match where
# the author last left the cursor.

$wb = $excel.ActiveWorkbook

# --- "General" sheet: append the new backlog row --------------------------
$ws1 = $wb.Worksheets.Item("General")

$ws1.Range("A19").Value = 18
$ws1.Range("B19").Value = "Escribir un cliente perfecto"
$ws1.Range("C19").Value = "Que conteste bien todas las preguntas y tenga puntaje perfecto."

# --- "Agregar jugadores" sheet: move the saved selection to D4 ------------
$ws2 = $wb.Worksheets.Item("Agregar jugadores")
[void]$ws2.Range("D4").Select()

# --- Back on "General": scroll so row 4 is at the top, keep D18 selected,
#     and leave this as the active sheet/tab (matches the source file). ---
[void]$ws1.Activate()
[void]$ws1.Range("D18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
